# Tab_2a_Bereiche: insert a new row for "Z04_B03 / Soziale Lage und Bildung"
# between the existing "Z04_B02" row (row 10) and the "Z05_B01" row
# (old row 11), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11, pushing rows 11..41 down to 12..42.
$ws.Rows.Item(11).Insert()

# The freshly inserted row picks up a generic default style; copy the
# formatting from the row above (row 10) so the new row matches the rest
# of the table (style index 4 in the original workbook).
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row's values.
$ws.Range("A11").Value2 = "Z04_B03"
$ws.Range("B11").Value2 = "Z04"
$ws.Range("C11").Value2 = "Soziale Lage und Bildung"
$ws.Range("D11").Value2 = "XXXSoziale Lage und Bildung"
